# "catégorie ok manque stocke available"
#  - L2 (Cétégorie N°ID) : 557 -> 558 (catégorie fixée)
#  - AD2 (Etat N°ID) : 1 -> "new" (nouvelle valeur texte "stock available" manquante)
#  - la sélection / vue se déplace vers la cellule qui vient d'être éditée (L2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cétégorie N°ID : incrémente l'identifiant de catégorie
$ws.Range("L2").Value = 558

# Etat (neuf1 utilisé2 reconditionné3) N°ID : nouvelle valeur texte "new"
$ws.Range("AD2").Value = "new"

# Ramène la vue / sélection sur L2 (comme dans le classeur édité)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("L2").Select()
